$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.110176086425781
$ws.Range("B1").Value = 1.553109407424927
$ws.Range("C1").Value = 3.985412836074829
$ws.Range("D1").Value = 1.500030159950256
$ws.Range("E1").Value = 0.9808575510978699
